$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Last status check on" timestamp in F1 (17:15 -> 17:30)
$ws.Range("F1").Value = "Last status check on: 03.02.2022 17:30"

# 2. D7: convert from text "+0.6" to numeric 0.6
$ws.Range("D7").Value = 0.6

# 3. E7: convert from text "2022-02-03 17:15:24" to numeric date serial,
#    with the same date/time number format style used by the other rows
#    in column E (copy style from E6, which already carries the date format).
$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("E7").Value = Get-Date -Year 2022 -Month 2 -Day 3 -Hour 17 -Minute 15 -Second 24
